# update backlog sprint 2 day 3
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sprint 2")

# Day 3 (column F) hours for the first three backlog rows
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 0

# Make "Sprint 2" the active sheet/tab with the given selection,
# matching the saved workbook + sheet view state.
$ws.Activate()
$ws.Range("F7").Select()
